# Generate Report for Handoff
# - Update "Status" text from "Handed back: in sync with en-US" to "Ready for handoff"
#   on the Overview sheet (zh-cn/de-de status columns) and on the per-locale
#   "zh-cn" / "de-de" sheets (Status column).
# - Bump the "Latest HO Xliff Generate Date" / "Latest Handback DateTime" timestamp
#   on the Overview and de-de sheets.
# - Bump the "Latest Handoff Datetime" timestamp on the zh-cn sheet.
# - Shrink the now-narrower Status columns to fit the shorter text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# Overview sheet: zh-cn / de-de status cells
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
# Overview sheet: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-27 06:56:15"

# zh-cn sheet: Status cell
$wsZhCn.Range("C2").Value = $newStatus
# zh-cn sheet: Latest Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-08-27 06:56:11"

# de-de sheet: Status cell
$wsDeDe.Range("C2").Value = $newStatus
# de-de sheet: Latest Handback DateTime
$wsDeDe.Range("H2").Value = "2016-08-27 06:56:15"

# Narrow the Status columns now that the text is shorter.
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
